$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 500.1
$ws.Range("J17").Value = 500.1
$ws.Range("L17").Value = 1500.3
$ws.Range("N17").Value = -1836.3
$ws.Range("H125").Value = 1662.8
$ws.Range("I125").Value = 1578.4
$ws.Range("K125").Value = 14205.6
$ws.Range("M125").Value = -11745.6
$ws.Range("H138").Value = 1828.3765
$ws.Range("I138").Value = 1361.0588
$ws.Range("J138").Value = 1945.2059
$ws.Range("K138").Value = 4083.1764
$ws.Range("L138").Value = 5835.6177
$ws.Range("M138").Value = 1056.8236
$ws.Range("N138").Value = -16115.6177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2750.052
$ws.Range("I32").Value = 2760.25
$ws.Range("J32").Value = 2673
$ws.Range("K32").Value = 2760.25
$ws.Range("L32").Value = 2673
$ws.Range("M32").Value = -2473.25
$ws.Range("N32").Value = -3247
$ws.Range("H61").Value = 1068.5416
$ws.Range("I61").Value = 1064.5238
$ws.Range("J61").Value = 1096.6666
$ws.Range("K61").Value = 1064.5238
$ws.Range("L61").Value = 1096.6666
$ws.Range("M61").Value = -852.5237999999999
$ws.Range("N61").Value = -1520.6666
$ws.Range("H74").Value = 1155.1
$ws.Range("I74").Value = 760.5789
$ws.Range("K74").Value = 760.5789
$ws.Range("M74").Value = 113.4211
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
$ws.Range("H77").Value = 1155.1
$ws.Range("I77").Value = 760.5789
$ws.Range("K77").Value = 3802.8945
$ws.Range("M77").Value = 565.1055000000001
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
$ws.Range("H88").Value = 2744.3333
$ws.Range("I88").Value = 2266.3333
$ws.Range("J88").Value = 2983.3333
$ws.Range("K88").Value = 2266.3333
$ws.Range("L88").Value = 2983.3333
$ws.Range("M88").Value = -1860.3333
$ws.Range("N88").Value = -3795.3333
$ws.Range("H91").Value = 2744.3333
$ws.Range("I91").Value = 2266.3333
$ws.Range("J91").Value = 2983.3333
$ws.Range("K91").Value = 2266.3333
$ws.Range("L91").Value = 2983.3333
$ws.Range("M91").Value = -862.3332999999998
$ws.Range("N91").Value = -5791.3333
$ws.Range("H110").Value = 1439.4117
$ws.Range("I110").Value = 1003.6667
$ws.Range("K110").Value = 1003.6667
$ws.Range("M110").Value = 1041.3333
$ws.Range("H114").Value = 21554.889
$ws.Range("J114").Value = 21554.889
$ws.Range("L114").Value = 21554.889
$ws.Range("N114").Value = -30232.889
$ws.Range("H122").Value = 1214
$ws.Range("I122").Value = 1033
$ws.Range("K122").Value = 3099
$ws.Range("M122").Value = -649
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null
$ws.Range("H132").Value = 2061.2222
$ws.Range("J132").Value = 3999.5
$ws.Range("L132").Value = 11998.5
$ws.Range("N132").Value = -17058.5
$ws.Range("H136").Value = 1068.5416
$ws.Range("I136").Value = 1064.5238
$ws.Range("J136").Value = 1096.6666
$ws.Range("K136").Value = 3193.5714
$ws.Range("L136").Value = 3289.9998
$ws.Range("M136").Value = -643.5713999999998
$ws.Range("N136").Value = -8389.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4124.48
$ws.Range("I86").Value = 4195.857
$ws.Range("K86").Value = 4195.857
$ws.Range("M86").Value = -3072.857
$ws.Range("H89").Value = 4124.48
$ws.Range("I89").Value = 4195.857
$ws.Range("K89").Value = 20979.285
$ws.Range("M89").Value = -15363.285
$ws.Range("H134").Value = 5321.1377
$ws.Range("I134").Value = 1212.2778
$ws.Range("K134").Value = 3636.8334
$ws.Range("M134").Value = -1101.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1191.2363
$ws.Range("I31").Value = 1167
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1167
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -872
$ws.Range("N31").Value = -3090
$ws.Range("H34").Value = 1191.2363
$ws.Range("I34").Value = 1167
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1167
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -965
$ws.Range("N34").Value = -2904
$ws.Range("H132").Value = 1845.24
$ws.Range("I132").Value = 1144.9375
$ws.Range("J132").Value = 3090.2222
$ws.Range("K132").Value = 3434.8125
$ws.Range("L132").Value = 9270.6666
$ws.Range("M132").Value = -904.8125
$ws.Range("N132").Value = -14330.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1122.7142
$ws.Range("I11").Value = 1261.5
$ws.Range("K11").Value = 3784.5
$ws.Range("M11").Value = -3644.5
$ws.Range("H12").Value = 71.23077000000001
$ws.Range("J12").Value = 62.8
$ws.Range("L12").Value = 188.4
$ws.Range("N12").Value = -534.4
$ws.Range("H57").Value = 1205
$ws.Range("I57").Value = 1205
$ws.Range("K57").Value = 3615
$ws.Range("M57").Value = -3056
$ws.Range("H58").Value = 2881
$ws.Range("J58").Value = 3100
$ws.Range("L58").Value = 9300
$ws.Range("N58").Value = -9556
$ws.Range("H62").Value = 6007
$ws.Range("J62").Value = 6007
$ws.Range("L62").Value = 18021
$ws.Range("N62").Value = -19393
$ws.Range("H63").Value = 6770.6665
$ws.Range("I63").Value = 6512
$ws.Range("J63").Value = 6900
$ws.Range("K63").Value = 19536
$ws.Range("L63").Value = 20700
$ws.Range("M63").Value = -18787
$ws.Range("N63").Value = -22198
$ws.Range("H64").Value = 4151.273
$ws.Range("I64").Value = 1980
$ws.Range("J64").Value = 4254.6665
$ws.Range("K64").Value = 5940
$ws.Range("L64").Value = 12763.9995
$ws.Range("M64").Value = -5670
$ws.Range("N64").Value = -13303.9995
$ws.Range("H65").Value = 6007
$ws.Range("J65").Value = 6007
$ws.Range("L65").Value = 54063
$ws.Range("N65").Value = -60927
$ws.Range("H66").Value = 6770.6665
$ws.Range("I66").Value = 6512
$ws.Range("J66").Value = 6900
$ws.Range("K66").Value = 58608
$ws.Range("L66").Value = 62100
$ws.Range("M66").Value = -54864
$ws.Range("N66").Value = -69588
$ws.Range("H67").Value = 4151.273
$ws.Range("I67").Value = 1980
$ws.Range("J67").Value = 4254.6665
$ws.Range("K67").Value = 5940
$ws.Range("L67").Value = 12763.9995
$ws.Range("M67").Value = -5004
$ws.Range("N67").Value = -14635.9995
$ws.Range("H107").Value = 11596.333
$ws.Range("J107").Value = 14739.571
$ws.Range("L107").Value = 44218.713
$ws.Range("N107").Value = -48058.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37503416
$ws.Range("I70").Value = 27781378
$ws.Range("K70").Value = 27781378
$ws.Range("M70").Value = -27781108
$ws.Range("H73").Value = 37503416
$ws.Range("I73").Value = 27781378
$ws.Range("K73").Value = 27781378
$ws.Range("M73").Value = -27780442
$ws.Range("H97").Value = 761.8
$ws.Range("I97").Value = 761.8
$ws.Range("K97").Value = 761.8
$ws.Range("M97").Value = -265.8
$ws.Range("H103").Value = 33332.332
$ws.Range("J103").Value = 33332.332
$ws.Range("L103").Value = 33332.332
$ws.Range("N103").Value = -35676.332
$ws.Range("H122").Value = 2599.7144
$ws.Range("I122").Value = 2699.6667
$ws.Range("K122").Value = 8099.000100000001
$ws.Range("M122").Value = -5649.000100000001
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H126").Value = 2047.2106
$ws.Range("I126").Value = 1716.6364
$ws.Range("J126").Value = 2501.75
$ws.Range("K126").Value = 5149.9092
$ws.Range("L126").Value = 7505.25
$ws.Range("M126").Value = -2679.9092
$ws.Range("N126").Value = -12445.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2985.5
$ws.Range("I40").Value = 2756.7144
$ws.Range("J40").Value = 3214.2856
$ws.Range("K40").Value = 2756.7144
$ws.Range("L40").Value = 3214.2856
$ws.Range("M40").Value = -2620.7144
$ws.Range("N40").Value = -3486.2856
$ws.Range("H69").Value = 35000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 35000
$ws.Range("M69").Value = $null
$ws.Range("N69").Value = -36622
$ws.Range("H72").Value = 35000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 105000
$ws.Range("M72").Value = $null
$ws.Range("N72").Value = -113112
$ws.Range("H122").Value = 35418788
$ws.Range("I122").Value = 56668570
$ws.Range("J122").Value = 2493.3333
$ws.Range("K122").Value = 170005710
$ws.Range("L122").Value = 7479.999899999999
$ws.Range("M122").Value = -170003260
$ws.Range("N122").Value = -12379.9999
$ws.Range("H136").Value = 1226.875
$ws.Range("I136").Value = 1178.3334
$ws.Range("J136").Value = 1566.6666
$ws.Range("K136").Value = 3535.0002
$ws.Range("L136").Value = 4699.9998
$ws.Range("M136").Value = -985.0001999999999
$ws.Range("N136").Value = -9799.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 23750
$ws.Range("J117").Value = 23750
$ws.Range("L117").Value = 23750
$ws.Range("N117").Value = -32928
$ws.Range("H122").Value = 57274136
$ws.Range("I122").Value = 74118936
$ws.Range("J122").Value = 1818
$ws.Range("K122").Value = 222356808
$ws.Range("L122").Value = 5454
$ws.Range("M122").Value = -222354358
$ws.Range("N122").Value = -10354
